$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "تحديد المستوى"
$ws.Range("A6").Value = "اختبار المستوى: اللفظي"
$ws.Range("A7").Value = "اختبار المستوى: الكمي"
$ws.Range("A8").Value = "الاختبارات"
$ws.Range("A9").Value = "الاختبار: اللفظي"
$ws.Range("A10").Value = "الاختبار: الكمي"
$ws.Range("A11").Value = "نصائح واستراتيجيات"
$ws.Range("A12").Value = "الإحصائيات"
$ws.Range("A13").Value = "نصمملك"

$ws.Range("B14").Select()
